# Add a new "Correction " column (N) to the Card11 sheet, right after the
# existing "Event" column (M), matching the data rows already present.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card11")

# Normalize the existing "Event " header text (drop trailing space).
$ws.Cells.Item(1, 13).Value = "Event"

# New header cell for the "Correction " column (keep trailing space) and
# copy the header style from the neighboring "Event" cell so formatting
# (bold, border, centered) matches the rest of the header row.
$ws.Cells.Item(1, 14).Value = "Correction "
$ws.Cells.Item(1, 14).Style = $ws.Cells.Item(1, 13).Style

# Fill data rows 2-12 in column M with "nan" (matching the sheet's
# convention for unset/empty values) and leave the new column N blank.
for ($row = 2; $row -le 12; $row++) {
    $ws.Cells.Item($row, 13).Value = "nan"
    $ws.Cells.Item($row, 14).Value = ""
}
